$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Group headers (merged cells A1/F1/K1) ---
$ws.Range("A1").Value = "Pittsburgh Pirates lefties"
$ws.Range("F1").Value = "Milwaukee Brewers righties"
$ws.Range("K1").Value = "Los Angeles Dodgers hitters"

# --- Group A: Pittsburgh Pirates lefties (columns A-D) ---
$ws.Range("A3").Value = "Frazier"
$ws.Range("B3").Value = 2200
$ws.Range("C3").Value = 0

$ws.Range("A4").Value = "Polanco"
$ws.Range("B4").Value = 3500
$ws.Range("C4").Formula = "=6.2+18.7"

$ws.Range("A5").Value = "Bell"
$ws.Range("B5").Value = 3200
$ws.Range("C5").Value = 9.7

$ws.Range("A6").Value = "Moran"
$ws.Range("B6").Value = 3100
$ws.Range("C6").Value = 3

# --- Group F: Milwaukee Brewers righties (columns F-I) ---
$ws.Range("F3").Value = "Cain"
$ws.Range("G3").Value = 4300
$ws.Range("H3").Value = 3

$ws.Range("F4").Value = "Braun"
$ws.Range("G4").Value = 3500
$ws.Range("H4").Value = 3

$ws.Range("F5").Value = "Santana"
$ws.Range("G5").Value = 3100
$ws.Range("H5").Value = 3

$ws.Range("F6").Value = "Perez"
$ws.Range("G6").Value = 3200
$ws.Range("H6").Value = 3

# --- Group K: Los Angeles Dodgers hitters (columns K-N) ---
$ws.Range("K3").Value = "Utley"
$ws.Range("L3").Value = 2500
$ws.Range("M3").Value = 0

$ws.Range("K4").Value = "Pederson"
$ws.Range("L4").Value = 2500
$ws.Range("M4").Value = 9.5

$ws.Range("K5").Value = "Bellinger"
$ws.Range("L5").Value = 3300
$ws.Range("M5").Value = 18.7

$ws.Range("K6").Value = "Muncy"
$ws.Range("L6").Value = 2200
$ws.Range("M6").Value = 9.5

# --- Success/Failure labels under each stack (row 8) ---
$ws.Range("D8").ClearContents()
$ws.Range("I8").Value = "Failure"
$ws.Range("N8").ClearContents()

# --- Update the remembered selection to match the authored edit ---
$ws.Activate()
$ws.Range("K9").Select()
